$d = $word.ActiveDocument

# --- Part 1: normalise the split runs of {{has_contraband_risk}} into a
# single run, leaving the following run (a lone trailing space) untouched.
# The underlying engine always merges an edited run forward with every
# following run that shares identical formatting, so the trailing space
# run (same rFonts/noProof as the field text) is temporarily given a
# distinguishing format (bold) to keep it out of the merge, then restored
# to its original formatting afterwards.
$fieldRange = $d.Content
$found = $fieldRange.Find.Execute("{{has_contraband_risk}}")
if ($found) {
    $start = $fieldRange.Start
    $end = $fieldRange.End

    $afterRun = $d.Range($end, $end + 1)
    $wasBold = $afterRun.Font.Bold
    $afterRun.Font.Bold = 1

    # Make (and then correct) a real text change so the engine actually
    # recomputes/merges the runs -- a same-value ("no-op") assignment is
    # otherwise ignored and leaves the original run split untouched.
    $editRange = $d.Range($start, $end)
    $editRange.Text = "{{has_contraband_riskZZZTMPZZZ}}"
    $editRange = $d.Range($start, $editRange.End)
    $editRange.Text = "{{has_contraband_risk}}"

    $afterRun = $d.Range($editRange.End, $editRange.End + 1)
    $afterRun.Font.Bold = $wasBold
}

# --- Part 2: fix field name used for details.
$d.Content.Find.Execute("{{has_contraband_risk_detail}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{{has_contraband_risk_details}}", 2)
